$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "old_site" column (B) values were a placeholder "Test Site" for every
# row. Unassigned APs should instead have an empty old_site, so clear the
# data cells in column B (rows 2-4) while leaving the B1 header ("old_site")
# intact.
$ws.Range("B2:B4").ClearContents()

# Update the active selection to match (was C2:C4 / active C2).
$ws.Range("B2:B4").Select() | Out-Null
